$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting all existing
# columns (A:V) one position to the right (-> B:W).
$ws.Columns("A:A").Insert()

# New column A becomes "Match ID": header label in row 2, and the
# match id value (7) for every visible data row (4-19).
$ws.Range("A2").Value = "Match ID"
$ws.Range("A4:A19").Value = 7

# Row 20 is a hidden "totals" row. Writing straight into a hidden row's
# blank cell causes the row to pick up a stray explicit height, so
# temporarily unhide it, set the value, then hide it again.
$ws.Rows(20).Hidden = $false
$ws.Cells.Item(20, 1).Value = 7
$ws.Rows(20).Hidden = $true

# Bold, borderless style (like the other header/data cells use for
# their font, but without the thin-border) for the new Match ID column
# header + visible data cells.
$ws.Range("A2:A19").Font.Bold = $true

# Restore a sensible selection over the newly added column's data.
$null = $ws.Range("A2:A19").Select()
